# Fruta / hortaliza, semanal
# The weekly refresh re-sorted the data rows (2-16); each destination row's
# full contents (A:T) come from a specific source row in the original sheet.
# Snapshot every source row first (so overwrites don't clobber a value we
# still need), then write them all into their destination rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (values that must land in destinationRow are
# the original contents of sourceRow)
$mapping = @{
    2  = 15
    3  = 4
    4  = 5
    5  = 2
    6  = 7
    7  = 8
    8  = 13
    9  = 16
    10 = 11
    11 = 6
    12 = 3
    13 = 12
    14 = 14
    15 = 9
    16 = 10
}

# Snapshot original row contents (columns A:T) before any writes happen.
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = $ws.Range("A$srcRow`:T$srcRow").Value2
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $ws.Range("A$destRow`:T$destRow").Value2 = $snapshot[$srcRow]
}
